$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 86
$ws.Range("H86").Value = 2748.75
$ws.Range("I86").Value = 2331.6667
$ws.Range("K86").Value = 2331.6667
$ws.Range("M86").Value = -1208.6667
# row 89
$ws.Range("H89").Value = 2748.75
$ws.Range("I89").Value = 2331.6667
$ws.Range("K89").Value = 11658.3335
$ws.Range("M89").Value = -6042.333500000001
# row 95
$ws.Range("H95").Value = 64000
$ws.Range("J95").Value = 64000
$ws.Range("L95").Value = 64000
$ws.Range("N95").Value = -69492
# row 137
$ws.Range("H137").Value = 31262864
$ws.Range("I137").Value = 125000750
$ws.Range("K137").Value = 375002250
$ws.Range("M137").Value = -374999700
# row 138
$ws.Range("H138").Value = 2346.8865
$ws.Range("I138").Value = 1541.7368
$ws.Range("J138").Value = 2958.8
$ws.Range("K138").Value = 4625.2104
$ws.Range("L138").Value = 8876.400000000001
$ws.Range("M138").Value = 514.7896000000001
$ws.Range("N138").Value = -19156.4
# row 141
$ws.Range("H141").Value = 1100
$ws.Range("I141").Value = 1100
$ws.Range("K141").Value = 3300
$ws.Range("M141").Value = 1880

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 7830.1177
$ws.Range("I61").Value = 7042.364
$ws.Range("J61").Value = 9274.333000000001
$ws.Range("K61").Value = 7042.364
$ws.Range("L61").Value = 9274.333000000001
$ws.Range("M61").Value = -6830.364
$ws.Range("N61").Value = -9698.333000000001
# row 63
$ws.Range("H63").Value = 7638.1665
$ws.Range("I63").Value = 2222.2856
$ws.Range("J63").Value = 9286.478999999999
$ws.Range("K63").Value = 2222.2856
$ws.Range("L63").Value = 9286.478999999999
$ws.Range("M63").Value = -1536.2856
$ws.Range("N63").Value = -10658.479
# row 66
$ws.Range("H66").Value = 7638.1665
$ws.Range("I66").Value = 2222.2856
$ws.Range("J66").Value = 9286.478999999999
$ws.Range("K66").Value = 11111.428
$ws.Range("L66").Value = 46432.395
$ws.Range("M66").Value = -7679.428
$ws.Range("N66").Value = -53296.395
# row 69
$ws.Range("H69").Value = 566666.7
$ws.Range("J69").Value = 566666.7
$ws.Range("L69").Value = 566666.7
$ws.Range("N69").Value = -568164.7
# row 72
$ws.Range("H72").Value = 566666.7
$ws.Range("J72").Value = 566666.7
$ws.Range("L72").Value = 1700000.1
$ws.Range("N72").Value = -1707488.1
# row 97
$ws.Range("H97").Value = 2058483
$ws.Range("I97").Value = 2849915.2
$ws.Range("K97").Value = 2849915.2
$ws.Range("M97").Value = -2849419.2
# row 122
$ws.Range("H122").Value = 2777.8333
$ws.Range("I122").Value = 1911
$ws.Range("K122").Value = 5733
$ws.Range("M122").Value = -3283
# row 132
$ws.Range("H132").Value = 4601.316
$ws.Range("I132").Value = 3231
$ws.Range("K132").Value = 9693
$ws.Range("M132").Value = -7163
# row 136
$ws.Range("H136").Value = 7830.1177
$ws.Range("I136").Value = 7042.364
$ws.Range("J136").Value = 9274.333000000001
$ws.Range("K136").Value = 21127.092
$ws.Range("L136").Value = 27822.999
$ws.Range("M136").Value = -18577.092
$ws.Range("N136").Value = -32922.999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Range("H94").Value = 1624.6487
$ws.Range("I94").Value = 1727.2084
$ws.Range("K94").Value = 1727.2084
$ws.Range("M94").Value = -1276.2084
# row 132
$ws.Range("H132").Value = 119887
$ws.Range("J132").Value = 119887
$ws.Range("L132").Value = 119887
$ws.Range("N132").Value = -130007
# row 134
$ws.Range("H134").Value = 20088.2
$ws.Range("I134").Value = 32970.75
$ws.Range("J134").Value = 11499.833
$ws.Range("K134").Value = 98912.25
$ws.Range("L134").Value = 34499.499
$ws.Range("M134").Value = -96377.25
$ws.Range("N134").Value = -39569.499

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 3
$ws.Range("H3").Value = 2264.111
$ws.Range("I3").Value = 672.125
$ws.Range("K3").Value = 672.125
$ws.Range("M3").Value = -559.125
# row 31
$ws.Range("H31").Value = 58828628
$ws.Range("I31").Value = 200001100
$ws.Range("K31").Value = 200001100
$ws.Range("M31").Value = -200000805
# row 34
$ws.Range("H34").Value = 58828628
$ws.Range("I34").Value = 200001100
$ws.Range("K34").Value = 200001100
$ws.Range("M34").Value = -200000898
# row 59
$ws.Range("H59").Value = 55000
$ws.Range("J59").Value = 65000
$ws.Range("L59").Value = 65000
$ws.Range("N59").Value = -67290

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 35
$ws.Range("H35").Value = 3315.5715
$ws.Range("I35").Value = 433
$ws.Range("J35").Value = 5477.5
$ws.Range("K35").Value = 1299
$ws.Range("L35").Value = 16432.5
$ws.Range("M35").Value = -1011
$ws.Range("N35").Value = -17008.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 3
$ws.Range("H3").Value = 286276
$ws.Range("I3").Value = 333929.66
$ws.Range("J3").Value = 354
$ws.Range("K3").Value = 333929.66
$ws.Range("L3").Value = 354
$ws.Range("M3").Value = -333813.66
$ws.Range("N3").Value = -586
# row 10
$ws.Range("H10").Value = 42647.9
$ws.Range("J10").Value = 47275.445
$ws.Range("L10").Value = 47275.445
$ws.Range("N10").Value = -47613.445
# row 11
$ws.Range("H11").Value = 3112353.5
$ws.Range("I11").Value = 231839.31
$ws.Range("J11").Value = 10601691
$ws.Range("K11").Value = 231839.31
$ws.Range("L11").Value = 10601691
$ws.Range("M11").Value = -231700.31
$ws.Range("N11").Value = -10601969
# row 12
$ws.Range("H12").Value = 269500
$ws.Range("J12").Value = 8500
$ws.Range("L12").Value = 8500
$ws.Range("N12").Value = -8780
# row 14
$ws.Range("H14").Value = 806294.25
$ws.Range("I14").Value = 750235.4399999999
$ws.Range("J14").Value = 1002500
$ws.Range("K14").Value = 750235.4399999999
$ws.Range("L14").Value = 1002500
$ws.Range("M14").Value = -750067.4399999999
$ws.Range("N14").Value = -1002836
# row 80
$ws.Range("H80").Value = 2445.1177
$ws.Range("I80").Value = 2371.9092
$ws.Range("J80").Value = 2579.3333
$ws.Range("K80").Value = 2371.9092
$ws.Range("L80").Value = 2579.3333
$ws.Range("M80").Value = -1373.9092
$ws.Range("N80").Value = -4575.3333
# row 83
$ws.Range("H83").Value = 2445.1177
$ws.Range("I83").Value = 2371.9092
$ws.Range("J83").Value = 2579.3333
$ws.Range("K83").Value = 11859.546
$ws.Range("L83").Value = 12896.6665
$ws.Range("M83").Value = -6867.546
$ws.Range("N83").Value = -22880.6665
# row 97
$ws.Range("H97").Value = 1300.3334
$ws.Range("I97").Value = 1534.6666
$ws.Range("K97").Value = 1534.6666
$ws.Range("M97").Value = -1038.6666
# row 122
$ws.Range("H122").Value = 6375.3125
$ws.Range("I122").Value = 8570
$ws.Range("J122").Value = 4180.625
$ws.Range("K122").Value = 25710
$ws.Range("L122").Value = 12541.875
$ws.Range("M122").Value = -23260
$ws.Range("N122").Value = -17441.875
# row 132
$ws.Range("H132").Value = 7749.3335
$ws.Range("I132").Value = 3986.625
$ws.Range("K132").Value = 11959.875
$ws.Range("M132").Value = -9429.875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# row 40
$ws.Range("H40").Value = 3235.0715
$ws.Range("I40").Value = 2208.2727
$ws.Range("K40").Value = 2208.2727
$ws.Range("M40").Value = -2072.2727
# row 46
$ws.Range("H46").Value = 9550.071
$ws.Range("I46").Value = 2999
$ws.Range("J46").Value = 10054
$ws.Range("K46").Value = 2999
$ws.Range("L46").Value = 10054
$ws.Range("M46").Value = -2811
$ws.Range("N46").Value = -10430
# row 82
$ws.Range("H82").Value = 3307.7273
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 3798.125
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 3798.125
$ws.Range("M82").Value = -1639
$ws.Range("N82").Value = -4520.125
# row 85
$ws.Range("H85").Value = 3307.7273
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 3798.125
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 3798.125
$ws.Range("M85").Value = -752
$ws.Range("N85").Value = -6294.125
# row 122
$ws.Range("H122").Value = 4377.9
$ws.Range("I122").Value = 3824.8572
$ws.Range("J122").Value = 5668.3335
$ws.Range("K122").Value = 11474.5716
$ws.Range("L122").Value = 17005.0005
$ws.Range("M122").Value = -9024.571599999999
$ws.Range("N122").Value = -21905.0005
# row 132
$ws.Range("H132").Value = 10737.435
$ws.Range("I132").Value = 10591.4375
$ws.Range("J132").Value = 11071.143
$ws.Range("K132").Value = 31774.3125
$ws.Range("L132").Value = 33213.429
$ws.Range("M132").Value = -29244.3125
$ws.Range("N132").Value = -38273.429

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 5
$ws.Range("H5").Value = 17666.666
$ws.Range("J5").Value = 17666.666
$ws.Range("L5").Value = 17666.666
$ws.Range("N5").Value = -17890.666
# row 12
$ws.Range("H12").Value = 99
$ws.Range("I12").Value = 99
$ws.Range("K12").Value = 99
$ws.Range("M12").Value = 43
# row 81
$ws.Range("H81").Value = 11804.892
$ws.Range("I81").Value = 1555.7858
$ws.Range("J81").Value = 18043.479
$ws.Range("K81").Value = 3111.5716
$ws.Range("L81").Value = 36086.958
$ws.Range("M81").Value = -2050.5716
$ws.Range("N81").Value = -38208.958
# row 84
$ws.Range("H84").Value = 11804.892
$ws.Range("I84").Value = 1555.7858
$ws.Range("J84").Value = 18043.479
$ws.Range("K84").Value = 15557.858
$ws.Range("L84").Value = 180434.79
$ws.Range("M84").Value = -10253.858
$ws.Range("N84").Value = -191042.79
# row 122
$ws.Range("H122").Value = 4609
$ws.Range("I122").Value = 4667.3687
$ws.Range("K122").Value = 14002.1061
$ws.Range("M122").Value = -11552.1061
# row 123
$ws.Range("H123").Value = 64500
$ws.Range("J123").Value = 64500
$ws.Range("L123").Value = 64500
$ws.Range("N123").Value = -74300
# row 132
$ws.Range("H132").Value = 4954
$ws.Range("I132").Value = 1804.1666
$ws.Range("J132").Value = 9315.308000000001
$ws.Range("K132").Value = 5412.4998
$ws.Range("L132").Value = 27945.924
$ws.Range("M132").Value = -2882.4998
$ws.Range("N132").Value = -33005.924

